$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "“XGBoost와 사이킷런을 활용한 그레이디언트 부스팅”이 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/04/06/xgboost%ec%99%80-%ec%82%ac%ec%9d%b4%ed%82%b7%eb%9f%b0%ec%9d%84-%ed%99%9c%ec%9a%a9%ed%95%9c-%ea%b7%b8%eb%a0%88%ec%9d%b4%eb%94%94%ec%96%b8%ed%8a%b8-%eb%b6%80%ec%8a%a4%ed%8c%85%ec%9d%b4-%ec%b6%9c/"

$ws.Range("D28").Value = "로봇 작업 할당 - Hungarian Algorithm"
$ws.Range("E28").Value = "https://ropiens.tistory.com/184"

$ws.Range("D37").Value = "[Paper Review]  N-HiTS: Neural Hierarchical Interpolation for Time Series Forecasting"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1974&mod=document&pageid=1"

$ws.Range("D46").Value = "대량 수혈 (massive transfusion)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/454"

$ws.Range("D50").Value = "incorrect theories"
$ws.Range("E50").Value = "http://incredible.egloos.com/7538902"
